# Apply updated cryptocurrency price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.490.01"
$ws.Range("E2").Value = "  +7.36%  "
$ws.Range("D3").Value = "3.583.40"
$ws.Range("E3").Value = "  +3.15%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'415.60"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").Value = "'129.10"
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("E7").Value = "  +3.51%  "
$ws.Range("D8").Value = "3.575.49"
$ws.Range("E8").Value = "  +3.11%  "
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("E10").Value = "  +5.82%  "
$ws.Range("E11").Value = "  +18.03%  "
$ws.Range("E12").Value = "  +47.03%  "
$ws.Range("D13").Value = "'42.31"
$ws.Range("E13").Value = "  -0.69%  "
$ws.Range("D14").Value = "'9.88"
$ws.Range("E14").Value = "  +1.64%  "
$ws.Range("D15").Value = "4.156.18"
$ws.Range("E15").Value = "  +3.33%  "
$ws.Range("E16").Value = "  -0.27%  "
$ws.Range("E17").Value = "  -0.96%  "
$ws.Range("D18").Value = "3.592.15"
$ws.Range("E18").Value = "  +3.27%  "
$ws.Range("E19").Value = "  +5.39%  "
$ws.Range("D20").Value = "67.355.06"
$ws.Range("E20").Value = "  +7.32%  "
$ws.Range("E21").Value = "  -3.06%  "
$ws.Range("D22").Value = "'451.22"
$ws.Range("E22").Value = "  -2.68%  "
$ws.Range("D23").Value = "'89.02"
$ws.Range("E23").Value = "  -1.85%  "
$ws.Range("D24").Value = "'3.14"
$ws.Range("E24").Value = "  -4.42%  "
$ws.Range("D25").Value = "'13.13"
$ws.Range("E25").Value = "  -1.19%  "
$ws.Range("E26").Value = "  +0.75%  "
$ws.Range("D27").Value = "'9.98"
$ws.Range("E27").Value = "  -7.30%  "
$ws.Range("D28").Value = "'34.99"
$ws.Range("E28").Value = "  +4.61%  "
$ws.Range("D29").Value = "'4.89"
$ws.Range("E29").Value = "  +1.91%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "'12.33"
$ws.Range("E30").Value = "  +1.97%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "'2.72"
$ws.Range("E31").Value = "  +1.88%  "
$ws.Range("E32").Value = "  +4.24%  "
$ws.Range("D33").Value = "'7.37"
$ws.Range("E33").Value = "  -2.56%  "
$ws.Range("D34").Value = "'0.161"
$ws.Range("E34").Value = "  -3.58%  "
$ws.Range("D35").Value = "'40.98"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Value = "'0.998"
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("E37").Value = "  -2.92%  "
$ws.Range("E38").Value = "  +0.38%  "
$ws.Range("D39").Value = "0.0₃0780"
$ws.Range("E39").Value = "  +36.86%  "
$ws.Range("E40").Value = "  +9.39%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "'0.998"
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'3.07"
$ws.Range("E42").Value = "  -0.99%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").Value = "'149.37"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "'2.75"
$ws.Range("E44").Value = "  +1.24%  "
$ws.Range("D45").Value = "'0.313"
$ws.Range("E45").Value = "  -2.64%  "
$ws.Range("D46").Value = "'3.25"
$ws.Range("E46").Value = "  -2.33%  "
$ws.Range("D47").Value = "'4.30"
$ws.Range("E47").Value = "  -2.08%  "
$ws.Range("E48").Value = "  -4.62%  "
$ws.Range("D49").Value = "'2.30"
$ws.Range("E49").Value = "  -3.30%  "
$ws.Range("D50").Value = "'115.65"
$ws.Range("E50").Value = "  +5.97%  "
$ws.Range("D51").Value = "'2.60"
$ws.Range("E51").Value = "  +11.28%  "
